$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 553 (the panda/snow post), shifting subsequent rows up.
$ws.Rows.Item(553).Delete()
